$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 1.103903333333333
$ws.Range("H2").Value = 3.31171
$ws.Range("I2").Value = 0.02393122995918198
$ws.Range("J2").Value = 0.02393122995918198
$ws.Range("M2").Value = 5.482938999999999
$ws.Range("N2").Value = 16.448817
$ws.Range("O2").Value = 0.1472261722051079
$ws.Range("P2").Value = 0.147226172205108
$ws.Range("Q2").Value = 6.052634638563333
$ws.Range("R2").Value = 54.47371174707
$ws.Range("S2").Value = 0.003523303383050565
$ws.Range("T2").Value = 0.003523303383050565

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 1.103903333333333
$ws.Range("H3").Value = 3.31171
$ws.Range("I3").Value = 0.02393122995918198
$ws.Range("J3").Value = 0.02393122995918198
$ws.Range("O3").Value = 0.5993885906243068
$ws.Range("P3").Value = 0.5993885906243068
$ws.Range("Q3").Value = 24.64154362797778
$ws.Range("R3").Value = 221.7738926518
$ws.Range("S3").Value = 0.01434410619714027
$ws.Range("T3").Value = 0.01434410619714027

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 1.103903333333333
$ws.Range("H4").Value = 3.31171
$ws.Range("I4").Value = 0.02393122995918198
$ws.Range("J4").Value = 0.02393122995918198
$ws.Range("M4").Value = 9.436472999999999
$ws.Range("N4").Value = 28.309419
$ws.Range("O4").Value = 0.2533852371705853
$ws.Range("P4").Value = 0.2533852371705853
$ws.Range("Q4").Value = 10.41695399961
$ws.Range("R4").Value = 93.75258599649
$ws.Range("S4").Value = 0.006063820378991142
$ws.Range("T4").Value = 0.006063820378991142

# Row 5 (FAPs -> ECs)
$ws.Range("H5").Value = 70.73212899999999
$ws.Range("I5").Value = 0.5111277390235027
$ws.Range("J5").Value = 0.5111277390235027
$ws.Range("M5").Value = 5.482938999999999
$ws.Range("N5").Value = 16.448817
$ws.Range("O5").Value = 0.1472261722051079
$ws.Range("P5").Value = 0.147226172205108
$ws.Range("Q5").Value = 129.2733162157103
$ws.Range("R5").Value = 1163.459845941393
$ws.Range("S5").Value = 0.07525138052428168
$ws.Range("T5").Value = 0.07525138052428169

# Row 6 (FAPs -> FAPs)
$ws.Range("H6").Value = 70.73212899999999
$ws.Range("I6").Value = 0.5111277390235027
$ws.Range("J6").Value = 0.5111277390235027
$ws.Range("O6").Value = 0.5993885906243068
$ws.Range("P6").Value = 0.5993885906243068
$ws.Range("Q6").Value = 526.2987528054243
$ws.Range("R6").Value = 4736.688775248818
$ws.Range("S6").Value = 0.3063641351222858
$ws.Range("T6").Value = 0.3063641351222858

# Row 7 (FAPs -> MuSCs)
$ws.Range("H7").Value = 70.73212899999999
$ws.Range("I7").Value = 0.5111277390235027
$ws.Range("J7").Value = 0.5111277390235027
$ws.Range("M7").Value = 9.436472999999999
$ws.Range("N7").Value = 28.309419
$ws.Range("O7").Value = 0.2533852371705853
$ws.Range("P7").Value = 0.2533852371705853
$ws.Range("Q7").Value = 222.487275180339
$ws.Range("R7").Value = 2002.38547662305
$ws.Range("S7").Value = 0.1295122233769353
$ws.Range("T7").Value = 0.1295122233769353

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 21.446869
$ws.Range("H8").Value = 64.34060699999999
$ws.Range("I8").Value = 0.4649410310173153
$ws.Range("J8").Value = 0.4649410310173154
$ws.Range("M8").Value = 5.482938999999999
$ws.Range("N8").Value = 16.448817
$ws.Range("O8").Value = 0.1472261722051079
$ws.Range("P8").Value = 0.147226172205108
$ws.Range("Q8").Value = 117.591874467991
$ws.Range("R8").Value = 1058.326870211919
$ws.Range("S8").Value = 0.0684514882977757
$ws.Range("T8").Value = 0.06845148829777573

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 21.446869
$ws.Range("H9").Value = 64.34060699999999
$ws.Range("I9").Value = 0.4649410310173153
$ws.Range("J9").Value = 0.4649410310173154
$ws.Range("O9").Value = 0.5993885906243068
$ws.Range("P9").Value = 0.5993885906243068
$ws.Range("Q9").Value = 478.7411562126732
$ws.Range("R9").Value = 4308.670405914059
$ws.Range("S9").Value = 0.2786803493048808
$ws.Range("T9").Value = 0.2786803493048808

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 21.446869
$ws.Range("H10").Value = 64.34060699999999
$ws.Range("I10").Value = 0.4649410310173153
$ws.Range("J10").Value = 0.4649410310173154
$ws.Range("M10").Value = 9.436472999999999
$ws.Range("N10").Value = 28.309419
$ws.Range("O10").Value = 0.2533852371705853
$ws.Range("P10").Value = 0.2533852371705853
$ws.Range("Q10").Value = 202.3828002530369
$ws.Range("R10").Value = 1821.445202277333
$ws.Range("S10").Value = 0.1178091934146589
$ws.Range("T10").Value = 0.1178091934146589
